$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume/1h (E) columns with refreshed crypto data.
# Numeric-looking Price strings need the cell pre-formatted as Text so
# Excel keeps them as literal strings (matching the source data which
# stores every Price/Volume cell as text, incl. values like "27.938.83").

$ws.Range('D2').Value = '27.929.34'
$ws.Range('E2').Value = '  +0.80%  '
$ws.Range('D3').Value = '1.880.54'
$ws.Range('E3').Value = '  +0.11%  '
$ws.Range('E4').Value = '  +1.67%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '335.23'
$ws.Range('E5').Value = '  +1.37%  '
$ws.Range('E6').Value = '  +1.52%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4696'
$ws.Range('E7').Value = '  -0.50%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3910'
$ws.Range('E8').Value = '  -1.48%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.14'
$ws.Range('E9').Value = '  -2.29%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07948'
$ws.Range('E10').Value = '  -1.13%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.008'
$ws.Range('E11').Value = '  -1.66%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.61'
$ws.Range('E12').Value = '  -1.14%  '
$ws.Range('D13').Value = '1.909.97'
$ws.Range('E13').Value = '  +1.60%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.934'
$ws.Range('E14').Value = '  -0.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.080'
$ws.Range('E15').Value = '  -1.38%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.020'
$ws.Range('E16').Value = '  +1.60%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.06762'
$ws.Range('E17').Value = '  +2.23%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.00001041'
$ws.Range('E19').Value = '  -0.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.03'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.018'
$ws.Range('D22').Value = '27.934.25'
$ws.Range('E22').Value = '  +0.79%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.471'
$ws.Range('E23').Value = '  -0.63%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.92'
$ws.Range('E24').Value = '  -0.98%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.363'
$ws.Range('E25').Value = '  +2.75%  '
$ws.Range('D26').Value = '2.132.03'
$ws.Range('E26').Value = '  +1.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '159.94'
$ws.Range('E27').Value = '  +2.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.88'
$ws.Range('E28').Value = '  -1.86%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.072'
$ws.Range('E29').Value = '  -1.11%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.452'
$ws.Range('E30').Value = '  -2.51%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '121.11'
$ws.Range('E31').Value = '  -1.21%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09506'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.9602'
$ws.Range('E33').Value = '  -1.33%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.662'
$ws.Range('E34').Value = '  +0.98%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.303'
$ws.Range('E35').Value = '  -0.17%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.344'
$ws.Range('E36').Value = '  -7.46%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06128'
$ws.Range('E37').Value = '  -0.08%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02236'
$ws.Range('E38').Value = '  -1.16%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.213'
$ws.Range('E39').Value = '  -1.22%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.017'
$ws.Range('E40').Value = '  +1.49%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.135'
$ws.Range('E41').Value = '  -0.26%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5914'
$ws.Range('E42').Value = '  -1.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1888'
$ws.Range('E43').Value = '  -0.88%  '
$ws.Range('E44').Value = '  -0.80%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.269'
$ws.Range('E45').Value = '  +0.73%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5636'
$ws.Range('E46').Value = '  -1.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.18'
$ws.Range('E47').Value = '  -0.81%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.397'
$ws.Range('E48').Value = '  -0.06%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.914'
$ws.Range('E49').Value = '  -1.07%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06881'
$ws.Range('E50').Value = '  +0.89%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '113.51'
$ws.Range('E51').Value = '  +1.66%  '
